$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-22 18:30:22"

for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
